$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 18 de Agosto de 2020 a las 13:04"

# Row 4
$ws.Cells.Item(4,2).Value = 5613183
$ws.Cells.Item(4,3).Value = 1156
$ws.Cells.Item(4,4).Value = 2974780
$ws.Cells.Item(4,5).Value = 2464631
$ws.Cells.Item(4,7).Value = 56
$ws.Cells.Item(4,8).Value = 173772

# Row 14
$ws.Cells.Item(14,2).Value = 347835
$ws.Cells.Item(14,3).Value = 2385
$ws.Cells.Item(14,4).Value = 300881
$ws.Cells.Item(14,5).Value = 26982
$ws.Cells.Item(14,7).Value = 168
$ws.Cells.Item(14,8).Value = 19972

# Row 38
$ws.Cells.Item(38,2).Value = 83418
$ws.Cells.Item(38,3).Value = 192
$ws.Cells.Item(38,4).Value = 77977
$ws.Cells.Item(38,5).Value = 4844
$ws.Cells.Item(38,7).Value = 9
$ws.Cells.Item(38,8).Value = 597

# Row 42
$ws.Cells.Item(42,2).Value = 72208
$ws.Cells.Item(42,3).Value = 1014
$ws.Cells.Item(42,4).Value = 33135
$ws.Cells.Item(42,5).Value = 35999
$ws.Cells.Item(42,7).Value = 45
$ws.Cells.Item(42,8).Value = 3074

# Row 59
$ws.Cells.Item(59,2).Value = 38449
$ws.Cells.Item(59,3).Value = 197
$ws.Cells.Item(59,5).Value = 3157
$ws.Cells.Item(59,7).Value = 1
$ws.Cells.Item(59,8).Value = 1992

# Row 69
$ws.Cells.Item(69,1).Value = "Nepal"
$ws.Cells.Item(69,2).Value = 28257
$ws.Cells.Item(69,3).Value = 1016
$ws.Cells.Item(69,4).Value = 17580
$ws.Cells.Item(69,5).Value = 10563
$ws.Cells.Item(69,7).Value = 7
$ws.Cells.Item(69,8).Value = 114

# Row 70
$ws.Cells.Item(70,1).Value = "Irlanda"
$ws.Cells.Item(70,2).Value = 27313
$ws.Cells.Item(70,4).Value = 23364
$ws.Cells.Item(70,5).Value = 2175
$ws.Cells.Item(70,8).Value = 1774

# Row 82
$ws.Cells.Item(82,2).Value = 14009
$ws.Cells.Item(82,3).Value = 123
$ws.Cells.Item(82,4).Value = 12767
$ws.Cells.Item(82,5).Value = 1069
$ws.Cells.Item(82,7).Value = 2
$ws.Cells.Item(82,8).Value = 173

# Row 84
$ws.Cells.Item(84,2).Value = 12485
$ws.Cells.Item(84,3).Value = 75
$ws.Cells.Item(84,4).Value = 6398
$ws.Cells.Item(84,5).Value = 5282
$ws.Cells.Item(84,7).Value = 2
$ws.Cells.Item(84,8).Value = 805

# Row 85
$ws.Cells.Item(85,2).Value = 12305
$ws.Cells.Item(85,3).Value = 68
$ws.Cells.Item(85,4).Value = 7767
$ws.Cells.Item(85,5).Value = 4282

# Row 91
$ws.Cells.Item(91,2).Value = 9219
$ws.Cells.Item(91,3).Value = 7
$ws.Cells.Item(91,4).Value = 8902
$ws.Cells.Item(91,5).Value = 192

# Row 98
$ws.Cells.Item(98,2).Value = 7776
$ws.Cells.Item(98,3).Value = 24
$ws.Cells.Item(98,5).Value = 392

# Row 144
$ws.Cells.Item(144,1).Value = "Malta"
$ws.Cells.Item(144,2).Value = 1423
$ws.Cells.Item(144,3).Value = 48
$ws.Cells.Item(144,4).Value = 766
$ws.Cells.Item(144,5).Value = 648
$ws.Cells.Item(144,8).Value = 9

# Row 145
$ws.Cells.Item(145,1).Value = "Jordania"
$ws.Cells.Item(145,2).Value = 1398
$ws.Cells.Item(145,4).Value = 1241
$ws.Cells.Item(145,5).Value = 146
$ws.Cells.Item(145,8).Value = 11

# Row 169
$ws.Cells.Item(169,2).Value = 486
$ws.Cells.Item(169,3).Value = 1
$ws.Cells.Item(169,5).Value = 29

# Row 184
$ws.Cells.Item(184,2).Value = 222
$ws.Cells.Item(184,3).Value = 5
$ws.Cells.Item(184,4).Value = 195
$ws.Cells.Item(184,5).Value = 27

# Row 213
$ws.Cells.Item(213,1).Value = "Islas Malvinas"
$ws.Cells.Item(213,4).Value = 13
$ws.Cells.Item(213,8).Value = 0

# Row 214
$ws.Cells.Item(214,1).Value = "Montserrat"
$ws.Cells.Item(214,4).Value = 12
$ws.Cells.Item(214,8).Value = 1

# Row 218
$ws.Cells.Item(218,2).Value = 5
$ws.Cells.Item(218,3).Value = 1
$ws.Cells.Item(218,5).Value = 4
